$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ntf3"
$ws.Range("C2").Value = "Ntrk1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.630712666666667
$ws.Range("H2").Value = 7.892138000000001
$ws.Range("I2").Value = 0.3947434022685045
$ws.Range("J2").Value = 0.3947434022685045
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.05351466666666666
$ws.Range("N2").Value = 0.160544
$ws.Range("O2").Value = 0.4186447970585551
$ws.Range("P2").Value = 0.4186447970585551
$ws.Range("Q2").Value = 0.1407817114524444
$ws.Range("R2").Value = 1.267035403072
$ws.Range("S2").Value = 0.1652572715329017
$ws.Range("T2").Value = 0.1652572715329017

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ntf3"
$ws.Range("C3").Value = "Ntrk1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.630712666666667
$ws.Range("H3").Value = 7.892138000000001
$ws.Range("I3").Value = 0.3947434022685045
$ws.Range("J3").Value = 0.3947434022685045
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.05507833333333334
$ws.Range("N3").Value = 0.165235
$ws.Range("O3").Value = 0.4308773485273219
$ws.Range("P3").Value = 0.4308773485273218
$ws.Range("Q3").Value = 0.1448952691588889
$ws.Range("R3").Value = 1.30405742243
$ws.Range("S3").Value = 0.1700859905181072
$ws.Range("T3").Value = 0.1700859905181072

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ntf3"
$ws.Range("C4").Value = "Ntrk1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.630712666666667
$ws.Range("H4").Value = 7.892138000000001
$ws.Range("I4").Value = 0.3947434022685045
$ws.Range("J4").Value = 0.3947434022685045
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.01923533333333333
$ws.Range("N4").Value = 0.057706
$ws.Range("O4").Value = 0.1504778544141231
$ws.Range("P4").Value = 0.1504778544141231
$ws.Range("Q4").Value = 0.05060263504755556
$ws.Range("R4").Value = 0.4554237154280001
$ws.Range("S4").Value = 0.05940014021749566
$ws.Range("T4").Value = 0.05940014021749565

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ntf3"
$ws.Range("C5").Value = "Ntrk1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.99311
$ws.Range("H5").Value = 11.97933
$ws.Range("I5").Value = 0.5991736942634763
$ws.Range("J5").Value = 0.5991736942634763
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.05351466666666666
$ws.Range("N5").Value = 0.160544
$ws.Range("O5").Value = 0.4186447970585551
$ws.Range("P5").Value = 0.4186447970585551
$ws.Range("Q5").Value = 0.2136899506133333
$ws.Range("R5").Value = 1.92320955552
$ws.Range("S5").Value = 0.2508409496377578
$ws.Range("T5").Value = 0.2508409496377578

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ntf3"
$ws.Range("C6").Value = "Ntrk1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.99311
$ws.Range("H6").Value = 11.97933
$ws.Range("I6").Value = 0.5991736942634763
$ws.Range("J6").Value = 0.5991736942634763
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.05507833333333334
$ws.Range("N6").Value = 0.165235
$ws.Range("O6").Value = 0.4308773485273219
$ws.Range("P6").Value = 0.4308773485273218
$ws.Range("Q6").Value = 0.2199338436166667
$ws.Range("R6").Value = 1.97940459255
$ws.Range("S6").Value = 0.2581703726915669
$ws.Range("T6").Value = 0.2581703726915669

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ntf3"
$ws.Range("C7").Value = "Ntrk1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.99311
$ws.Range("H7").Value = 11.97933
$ws.Range("I7").Value = 0.5991736942634763
$ws.Range("J7").Value = 0.5991736942634763
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.01923533333333333
$ws.Range("N7").Value = 0.057706
$ws.Range("O7").Value = 0.1504778544141231
$ws.Range("P7").Value = 0.1504778544141231
$ws.Range("Q7").Value = 0.07680880188666667
$ws.Range("R7").Value = 0.69127921698
$ws.Range("S7").Value = 0.0901623719341517
$ws.Range("T7").Value = 0.09016237193415169

$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Ntf3"
$ws.Range("C8").Value = "Ntrk1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.04053866666666667
$ws.Range("H8").Value = 0.121616
$ws.Range("I8").Value = 0.00608290346801924
$ws.Range("J8").Value = 0.006082903468019241
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.05351466666666666
$ws.Range("N8").Value = 0.160544
$ws.Range("O8").Value = 0.4186447970585551
$ws.Range("P8").Value = 0.4186447970585551
$ws.Range("Q8").Value = 0.002169413233777778
$ws.Range("R8").Value = 0.019524719104
$ws.Range("S8").Value = 0.002546575887895696
$ws.Range("T8").Value = 0.002546575887895696

$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Ntf3"
$ws.Range("C9").Value = "Ntrk1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.04053866666666667
$ws.Range("H9").Value = 0.121616
$ws.Range("I9").Value = 0.00608290346801924
$ws.Range("J9").Value = 0.006082903468019241
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.05507833333333334
$ws.Range("N9").Value = 0.165235
$ws.Range("O9").Value = 0.4308773485273219
$ws.Range("P9").Value = 0.4308773485273218
$ws.Range("Q9").Value = 0.002232802195555556
$ws.Range("R9").Value = 0.02009521976
$ws.Range("S9").Value = 0.002620985317647781
$ws.Range("T9").Value = 0.002620985317647781

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Ntf3"
$ws.Range("C10").Value = "Ntrk1"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.04053866666666667
$ws.Range("H10").Value = 0.121616
$ws.Range("I10").Value = 0.00608290346801924
$ws.Range("J10").Value = 0.006082903468019241
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.01923533333333333
$ws.Range("N10").Value = 0.057706
$ws.Range("O10").Value = 0.1504778544141231
$ws.Range("P10").Value = 0.1504778544141231
$ws.Range("Q10").Value = 0.0007797747662222223
$ws.Range("R10").Value = 0.007017972896
$ws.Range("S10").Value = 0.0009153422624757639
$ws.Range("T10").Value = 0.0009153422624757639
